$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2..519) from 45203 to 45204
$ws.Range("C2:C519").Value = 45204

# Row 519 picks up an explicit row height (matches the rows above it) once a new
# row is appended below it.
$ws.Rows.Item(519).RowHeight = 15

# Append new row 520 with the new "A 47553-2023" record
$ws.Range("A520").Value = "A 47553-2023"

$ws.Range("B520").NumberFormat = "YYYY-MM-DD"
$ws.Range("B520").Value = 45203

$ws.Range("C520").NumberFormat = "YYYY-MM-DD"
$ws.Range("C520").Value = 45204

$ws.Range("D520").Value = "UPPSALA LÄN"
$ws.Range("E520").Value = "ENKÖPING"

$ws.Range("G520").Value = 13.3
$ws.Range("H520").Value = 0
$ws.Range("I520").Value = 0
$ws.Range("J520").Value = 0
$ws.Range("K520").Value = 0
$ws.Range("L520").Value = 0
$ws.Range("M520").Value = 0
$ws.Range("N520").Value = 0
$ws.Range("O520").Value = 0
$ws.Range("P520").Value = 0
$ws.Range("Q520").Value = 0

$ws.Range("R520").WrapText = $true
